# Auto-generated Excel COM-interop script
# Commit: "Adicionados balancos concatenados em uma unica planilha."
# Adds columns V, W, X (periods 31/12/2023, 31/03/2024, 30/06/2024) to the
# MTRE3 consolidated balance-sheet worksheet, extending the used range from
# A1:U80 to A1:X80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy the style of U1 (bold, centered, bordered) into
#     V1:X1, then set the new period labels ---
$ws.Range("U1").Copy() | Out-Null
$ws.Range("V1:X1").PasteSpecial(-4122) | Out-Null
$ws.Range("V1").Value2 = "31/12/2023"
$ws.Range("W1").Value2 = "31/03/2024"
$ws.Range("X1").Value2 = "30/06/2024"

# --- Data rows 2-80: new quarterly figures in columns V (31/12/2023),
#     W (31/03/2024) and X (30/06/2024) ---
$ws.Cells.Item(2, 22).Value2 = 2165236.992
$ws.Cells.Item(2, 23).Value2 = 2059207.04
$ws.Cells.Item(2, 24).Value2 = 2051480.064
$ws.Cells.Item(3, 22).Value2 = 1658151.04
$ws.Cells.Item(3, 23).Value2 = 1600987.008
$ws.Cells.Item(3, 24).Value2 = 1742342.016
$ws.Cells.Item(4, 22).Value2 = 198786
$ws.Cells.Item(4, 23).Value2 = 188828.992
$ws.Cells.Item(4, 24).Value2 = 230530
$ws.Cells.Item(5, 22).Value2 = 0
$ws.Cells.Item(5, 23).Value2 = 0
$ws.Cells.Item(5, 24).Value2 = 0
$ws.Cells.Item(6, 22).Value2 = 599084.032
$ws.Cells.Item(6, 23).Value2 = 596035.968
$ws.Cells.Item(6, 24).Value2 = 665409.024
$ws.Cells.Item(7, 22).Value2 = 820078.976
$ws.Cells.Item(7, 23).Value2 = 741656
$ws.Cells.Item(7, 24).Value2 = 761732.992
$ws.Cells.Item(8, 22).Value2 = 0
$ws.Cells.Item(8, 23).Value2 = 0
$ws.Cells.Item(8, 24).Value2 = 0
$ws.Cells.Item(9, 22).Value2 = 19812
$ws.Cells.Item(9, 23).Value2 = 21601
$ws.Cells.Item(9, 24).Value2 = 22762
$ws.Cells.Item(10, 22).Value2 = 0
$ws.Cells.Item(10, 23).Value2 = 0
$ws.Cells.Item(10, 24).Value2 = 0
$ws.Cells.Item(11, 22).Value2 = 20390
$ws.Cells.Item(11, 23).Value2 = 52865
$ws.Cells.Item(11, 24).Value2 = 61908
$ws.Cells.Item(12, 22).Value2 = 418670.016
$ws.Cells.Item(12, 23).Value2 = 373636
$ws.Cells.Item(12, 24).Value2 = 227668.992
$ws.Cells.Item(13, 22).Value2 = 0
$ws.Cells.Item(13, 23).Value2 = 0
$ws.Cells.Item(13, 24).Value2 = 0
$ws.Cells.Item(14, 22).Value2 = 0
$ws.Cells.Item(14, 23).Value2 = 0
$ws.Cells.Item(14, 24).Value2 = 0
$ws.Cells.Item(15, 22).Value2 = 0
$ws.Cells.Item(15, 23).Value2 = 0
$ws.Cells.Item(15, 24).Value2 = 0
$ws.Cells.Item(16, 22).Value2 = 188743.008
$ws.Cells.Item(16, 23).Value2 = 213624.992
$ws.Cells.Item(16, 24).Value2 = 160302
$ws.Cells.Item(17, 22).Value2 = 176046
$ws.Cells.Item(17, 23).Value2 = 109035
$ws.Cells.Item(17, 24).Value2 = 16723
$ws.Cells.Item(18, 22).Value2 = 0
$ws.Cells.Item(18, 23).Value2 = 0
$ws.Cells.Item(18, 24).Value2 = 0
$ws.Cells.Item(19, 22).Value2 = 0
$ws.Cells.Item(19, 23).Value2 = 0
$ws.Cells.Item(19, 24).Value2 = 0
$ws.Cells.Item(20, 22).Value2 = 0
$ws.Cells.Item(20, 23).Value2 = 0
$ws.Cells.Item(20, 24).Value2 = 0
$ws.Cells.Item(21, 22).Value2 = 36144
$ws.Cells.Item(21, 23).Value2 = 33208
$ws.Cells.Item(21, 24).Value2 = 33209
$ws.Cells.Item(22, 22).Value2 = 30128
$ws.Cells.Item(22, 23).Value2 = 30599
$ws.Cells.Item(22, 24).Value2 = 32594
$ws.Cells.Item(23, 22).Value2 = 35814
$ws.Cells.Item(23, 23).Value2 = 32260
$ws.Cells.Item(23, 24).Value2 = 27861
$ws.Cells.Item(24, 22).Value2 = 22474
$ws.Cells.Item(24, 23).Value2 = 21725
$ws.Cells.Item(24, 24).Value2 = 21014
$ws.Cells.Item(25, 22).Value2 = 0
$ws.Cells.Item(25, 23).Value2 = 0
$ws.Cells.Item(25, 24).Value2 = 0
$ws.Cells.Item(26, 22).Value2 = 2165236.992
$ws.Cells.Item(26, 23).Value2 = 2059207.04
$ws.Cells.Item(26, 24).Value2 = 2051480.064
$ws.Cells.Item(27, 22).Value2 = 631276.992
$ws.Cells.Item(27, 23).Value2 = 613843.968
$ws.Cells.Item(27, 24).Value2 = 657304
$ws.Cells.Item(28, 22).Value2 = 20570
$ws.Cells.Item(28, 23).Value2 = 26697
$ws.Cells.Item(28, 24).Value2 = 21000
$ws.Cells.Item(29, 22).Value2 = 62398
$ws.Cells.Item(29, 23).Value2 = 42704
$ws.Cells.Item(29, 24).Value2 = 60052
$ws.Cells.Item(30, 22).Value2 = 21523
$ws.Cells.Item(30, 23).Value2 = 21739
$ws.Cells.Item(30, 24).Value2 = 24717
$ws.Cells.Item(31, 22).Value2 = 305188
$ws.Cells.Item(31, 23).Value2 = 297544.992
$ws.Cells.Item(31, 24).Value2 = 309750.016
$ws.Cells.Item(32, 22).Value2 = 0
$ws.Cells.Item(32, 23).Value2 = 0
$ws.Cells.Item(32, 24).Value2 = 0
$ws.Cells.Item(33, 22).Value2 = 0
$ws.Cells.Item(33, 23).Value2 = 0
$ws.Cells.Item(33, 24).Value2 = 0
$ws.Cells.Item(34, 22).Value2 = 219194
$ws.Cells.Item(34, 23).Value2 = 222516
$ws.Cells.Item(34, 24).Value2 = 239028
$ws.Cells.Item(35, 22).Value2 = 2404
$ws.Cells.Item(35, 23).Value2 = 2643
$ws.Cells.Item(35, 24).Value2 = 2757
$ws.Cells.Item(36, 22).Value2 = 0
$ws.Cells.Item(36, 23).Value2 = 0
$ws.Cells.Item(36, 24).Value2 = 0
$ws.Cells.Item(37, 22).Value2 = 507751.008
$ws.Cells.Item(37, 23).Value2 = 433188
$ws.Cells.Item(37, 24).Value2 = 372820
$ws.Cells.Item(38, 22).Value2 = 351479.008
$ws.Cells.Item(38, 23).Value2 = 293404.992
$ws.Cells.Item(38, 24).Value2 = 270227.008
$ws.Cells.Item(39, 22).Value2 = 0
$ws.Cells.Item(39, 23).Value2 = 0
$ws.Cells.Item(39, 24).Value2 = 0
$ws.Cells.Item(40, 22).Value2 = 131334
$ws.Cells.Item(40, 23).Value2 = 113313
$ws.Cells.Item(40, 24).Value2 = 78074
$ws.Cells.Item(41, 22).Value2 = 7606
$ws.Cells.Item(41, 23).Value2 = 8594
$ws.Cells.Item(41, 24).Value2 = 6416
$ws.Cells.Item(42, 22).Value2 = 0
$ws.Cells.Item(42, 23).Value2 = 0
$ws.Cells.Item(42, 24).Value2 = 0
$ws.Cells.Item(43, 22).Value2 = 17332
$ws.Cells.Item(43, 23).Value2 = 17876
$ws.Cells.Item(43, 24).Value2 = 18103
$ws.Cells.Item(44, 22).Value2 = 0
$ws.Cells.Item(44, 23).Value2 = 0
$ws.Cells.Item(44, 24).Value2 = 0
$ws.Cells.Item(45, 22).Value2 = 0
$ws.Cells.Item(45, 23).Value2 = 0
$ws.Cells.Item(45, 24).Value2 = 0
$ws.Cells.Item(46, 22).Value2 = 35802
$ws.Cells.Item(46, 23).Value2 = 35379
$ws.Cells.Item(46, 24).Value2 = 44236
$ws.Cells.Item(47, 22).Value2 = 990407.024
$ws.Cells.Item(47, 23).Value2 = 976795.976
$ws.Cells.Item(47, 24).Value2 = 977120.032
$ws.Cells.Item(48, 22).Value2 = 959492.992
$ws.Cells.Item(48, 23).Value2 = 959492.992
$ws.Cells.Item(48, 24).Value2 = 959492.992
$ws.Cells.Item(49, 22).Value2 = 0
$ws.Cells.Item(49, 23).Value2 = 0
$ws.Cells.Item(49, 24).Value2 = 0
$ws.Cells.Item(50, 22).Value2 = 0
$ws.Cells.Item(50, 23).Value2 = 0
$ws.Cells.Item(50, 24).Value2 = 0
$ws.Cells.Item(51, 22).Value2 = 30914
$ws.Cells.Item(51, 23).Value2 = 17303
$ws.Cells.Item(51, 24).Value2 = 17627
$ws.Cells.Item(52, 22).Value2 = 0
$ws.Cells.Item(52, 23).Value2 = 0
$ws.Cells.Item(52, 24).Value2 = 0
$ws.Cells.Item(53, 22).Value2 = 0
$ws.Cells.Item(53, 23).Value2 = 0
$ws.Cells.Item(53, 24).Value2 = 0
$ws.Cells.Item(54, 22).Value2 = 0
$ws.Cells.Item(54, 23).Value2 = 0
$ws.Cells.Item(54, 24).Value2 = 0
$ws.Cells.Item(55, 22).Value2 = 0
$ws.Cells.Item(55, 23).Value2 = 0
$ws.Cells.Item(55, 24).Value2 = 0
$ws.Cells.Item(56, 22).Value2 = 0
$ws.Cells.Item(56, 23).Value2 = 0
$ws.Cells.Item(56, 24).Value2 = 0
$ws.Cells.Item(57, 22).Value2 = "'"
$ws.Cells.Item(57, 23).Value2 = "'"
$ws.Cells.Item(57, 24).Value2 = "'"
$ws.Range("V57:X57").Style = "Normal"
$ws.Cells.Item(58, 22).Value2 = "'"
$ws.Cells.Item(58, 23).Value2 = "'"
$ws.Cells.Item(58, 24).Value2 = "'"
$ws.Range("V58:X58").Style = "Normal"
$ws.Cells.Item(59, 22).Value2 = 278712.992
$ws.Cells.Item(59, 23).Value2 = 245572.992
$ws.Cells.Item(59, 24).Value2 = 289632.992
$ws.Cells.Item(60, 22).Value2 = -224261.008
$ws.Cells.Item(60, 23).Value2 = -192352
$ws.Cells.Item(60, 24).Value2 = -220764.992
$ws.Cells.Item(61, 22).Value2 = 54451.992
$ws.Cells.Item(61, 23).Value2 = 53221
$ws.Cells.Item(61, 24).Value2 = 68868
$ws.Cells.Item(62, 22).Value2 = -18148
$ws.Cells.Item(62, 23).Value2 = -14741
$ws.Cells.Item(62, 24).Value2 = -15648
$ws.Cells.Item(63, 22).Value2 = -27499
$ws.Cells.Item(63, 23).Value2 = -26787
$ws.Cells.Item(63, 24).Value2 = -32368
$ws.Cells.Item(64, 22).Value2 = 0
$ws.Cells.Item(64, 23).Value2 = 0
$ws.Cells.Item(64, 24).Value2 = 0
$ws.Cells.Item(65, 22).Value2 = 12624
$ws.Cells.Item(65, 23).Value2 = 0
$ws.Cells.Item(65, 24).Value2 = 502
$ws.Cells.Item(66, 22).Value2 = 3087
$ws.Cells.Item(66, 23).Value2 = -6677
$ws.Cells.Item(66, 24).Value2 = 0
$ws.Cells.Item(67, 22).Value2 = 1095
$ws.Cells.Item(67, 23).Value2 = 520
$ws.Cells.Item(67, 24).Value2 = 2157
$ws.Cells.Item(68, 22).Value2 = -2090
$ws.Cells.Item(68, 23).Value2 = 2002
$ws.Cells.Item(68, 24).Value2 = 1943
$ws.Cells.Item(69, 22).Value2 = 5742
$ws.Cells.Item(69, 23).Value2 = 5905
$ws.Cells.Item(69, 24).Value2 = 5812
$ws.Cells.Item(70, 22).Value2 = -7832
$ws.Cells.Item(70, 23).Value2 = -3903
$ws.Cells.Item(70, 24).Value2 = -3869
$ws.Cells.Item(71, 22).Value2 = "'"
$ws.Cells.Item(71, 23).Value2 = "'"
$ws.Cells.Item(71, 24).Value2 = "'"
$ws.Range("V71:X71").Style = "Normal"
$ws.Cells.Item(72, 22).Value2 = "'"
$ws.Cells.Item(72, 23).Value2 = "'"
$ws.Cells.Item(72, 24).Value2 = "'"
$ws.Range("V72:X72").Style = "Normal"
$ws.Cells.Item(73, 22).Value2 = "'"
$ws.Cells.Item(73, 23).Value2 = "'"
$ws.Cells.Item(73, 24).Value2 = "'"
$ws.Range("V73:X73").Style = "Normal"
$ws.Cells.Item(74, 22).Value2 = 23521
$ws.Cells.Item(74, 23).Value2 = 7538
$ws.Cells.Item(74, 24).Value2 = 25454
$ws.Cells.Item(75, 22).Value2 = -4246
$ws.Cells.Item(75, 23).Value2 = -5056
$ws.Cells.Item(75, 24).Value2 = -6397
$ws.Cells.Item(76, 22).Value2 = -1698
$ws.Cells.Item(76, 23).Value2 = -428
$ws.Cells.Item(76, 24).Value2 = -210
$ws.Cells.Item(77, 22).Value2 = "'"
$ws.Cells.Item(77, 23).Value2 = "'"
$ws.Cells.Item(77, 24).Value2 = "'"
$ws.Range("V77:X77").Style = "Normal"
$ws.Cells.Item(78, 22).Value2 = "'"
$ws.Cells.Item(78, 23).Value2 = "'"
$ws.Cells.Item(78, 24).Value2 = "'"
$ws.Range("V78:X78").Style = "Normal"
$ws.Cells.Item(79, 22).Value2 = -5651
$ws.Cells.Item(79, 23).Value2 = -989
$ws.Cells.Item(79, 24).Value2 = -6523
$ws.Cells.Item(80, 22).Value2 = 13562
$ws.Cells.Item(80, 23).Value2 = 1065
$ws.Cells.Item(80, 24).Value2 = 12324

Write-Host "Done: added columns V:X (31/12/2023, 31/03/2024, 30/06/2024) for rows 1-80"
